# This script rewrites the body of the document (everything except the
# final sectPr, which Word preserves automatically) using Range.InsertXML
# so that paragraph/run/proofErr structure matches exactly what a human
# editor produced: dropped the numbered-list formatting on the first two
# bullet paragraphs (replaced by plain indentation), merged a couple of
# runs that ended up with identical formatting, relocated the "_GoBack"
# bookmark to the end of the document, tightened up the sentence
# "...drive electric. Tesla believes..." (with matching spelling/grammar
# proofErr markers), and appended a "Reference list:" section with the
# citation for the Telsa/Tesla official site.

$d = $word.ActiveDocument
$r = $d.Content

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:firstLineChars="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">established </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>brand</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>background</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/>
          <w:color w:val="222222"/>
          <w:szCs w:val="21"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Telsa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="222222"/>
          <w:szCs w:val="21"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> is an American automotive and energy company </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="666666"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>t</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>was</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> founded in 2003 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="666666"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">by a group of engineers who wanted to prove that people didn’t need to compromise to drive </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="666666"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>electric.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="666666"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Tesla</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="666666"/>
          <w:sz w:val="23"/>
          <w:szCs w:val="23"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> believes the faster the world stops relying on fossil fuels and move towards to clean energy ,the better future we would have;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
      </w:pPr>
      <w:r>
        <w:t>Reference list:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:t>elsa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>2008,</w:t>
      </w:r>
      <w:r>
        <w:t>Telsa</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>official site</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>&lt;</w:t>
      </w:r>
      <w:r>
        <w:t>tesla.com/about</w:t>
      </w:r>
      <w:r>
        <w:t>&gt;</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
